$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data values scraped on 2022-12-22 12:21 UTC
$ws.Range("D2").Value = "'246.04"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'22.70"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.408"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'0.05765"
$ws.Range("G5").Value = "'12"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'6.345"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'0.8104"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.8896"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.1453"
$ws.Range("E10").Value = "9WazirXWRXBestin24h"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.07329"
$ws.Range("G11").Value = "'12"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.02996"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "'12"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03114"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.09417"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'3.937"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.001577"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'0.04799"
$ws.Range("G17").Value = "'12"
$ws.Range("D18").Value = "'0.0005850"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'0.006408"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'0.004145"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.0009920"
$ws.Range("G21").Value = "'12"
$ws.Range("G22").Value = "'12"
$ws.Range("D23").Value = "'3.722"
$ws.Range("G23").Value = "'12"
$ws.Range("D24").Value = "'2.200"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.3271"
$ws.Range("G25").Value = "'12"
$ws.Range("G26").Value = "'12"
$ws.Range("D27").Value = "'0.0004650"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03907"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.006783"
$ws.Range("G41").Value = "'12"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.002420"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.006800"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005651"
$ws.Range("G45").Value = "'12"
$ws.Range("G46").Value = "'12"
$ws.Range("D47").Value = "'0.3800"
$ws.Range("G47").Value = "'12"
$ws.Range("D48").Value = "'0.1606"
$ws.Range("G48").Value = "'12"
$ws.Range("G49").Value = "'12"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
